# Apply the edits described by the diff to the "Journal de travail" workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# 1. Update the description text for the "Spécifications..." entry (row 8, column B)
#    - remove "mockup, " from the middle of the sentence
#    - append ", mockups" before the final period
$ws.Range("B8").Value = "Spécifications de notre projet, définition des souhaits de chacun, choix des outils à implémenter, schéma de Dépendances Fonctionnelles afin d'organiser l'ordre d'implémentation de chaque fonctionnalité et début de la rédaction du cahier des charges, mockups."

# 2. Update the hours worked for that same row (C8): 5 -> 6.5
$ws.Range("C8").Value = 6.5

# 3. Move the active cell selection to D8 (was B11)
$ws.Range("D8").Select()

$wb.Save()
